$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 44523.82425925926
$ws.Range("C5").Value = 44523.82746527778
$ws.Range("D5").Value = "IP Address"
$ws.Range("E5").Value = 100
$ws.Range("F5").Value = 277
$ws.Range("G5").Value = $true
$ws.Range("H5").Value = 44523.82747685185
$ws.Range("I5").Value = "3hostc"
$ws.Range("J5").Value = "ebola %>%`n  pivot_longer(``case_death``:last_col()) %>%`n  (separate, into = ____, sep = ____) %>%`n  drop_na()"

# Row 6
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 44523.82498842593
$ws.Range("C6").Value = 44523.82778935185
$ws.Range("D6").Value = "IP Address"
$ws.Range("E6").Value = 100
$ws.Range("F6").Value = 242
$ws.Range("G6").Value = $true
$ws.Range("H6").Value = 44523.82780092592
$ws.Range("I6").Value = "2dunic"
$ws.Range("J6").Value = "ebola_tidy <-ebola %>%`n  pivot_longer(289)____(____) %>%`n  ____(____, into = ____, sep = ____) %>%`n  ____"

# Row 7
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 44523.18524305556
$ws.Range("C7").Value = 44523.19250000001
$ws.Range("D7").Value = "Spam"
$ws.Range("E7").Value = 50
$ws.Range("F7").Value = 627
$ws.Range("G7").Value = $false
$ws.Range("H7").Value = 44523.87909722222
$ws.Range("I7").Value = "2nesch"

# Reset auto row heights so the multi-line q2 text in rows 5-6 doesn't
# leave an explicit customHeight behind (matches the plain <row> elements
# with no ht attribute in the target workbook).
$ws.Rows.Item(5).AutoFit()
$ws.Rows.Item(6).AutoFit()
